$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '68.481.01'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '3.808.23'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.35'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.24'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("D7").Value = '3.806.07'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.03'
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = '  +11.95%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.19'
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").Value = '4.446.79'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = '3.767.68'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").Value = '68.487.53'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.14'
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.09'
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '463.36'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.65'
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.03'
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.02'
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '3.958.46'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("E31").Value = '  -5.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.22'
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.26'
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.19'
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.05'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  +1.32%  '
$ws.Range("E38").Value = '  +8.26%  '
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.24'
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = '  -1.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.982'
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -0.57%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '153.33'
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.12'
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.10'
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = '  -2.98%  '
$ws.Range("E48").Value = '  -1.29%  '
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '381.90'
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = '  -1.75%  '
